$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Anodiam-Website")

$ws.Range("B7").Value = "End User"
$ws.Range("C7").Value = "open anodiam website"
$ws.Range("D7").Value = "I can see course catalogue"

$ws.Range("B8").Value = "End User"
$ws.Range("C8").Value = "Register as user"
$ws.Range("D8").Value = "I can do a course"

$ws.Range("B9").Value = "End User"
$ws.Range("C9").Value = "View the course details"
$ws.Range("D9").Value = "I can select the courses"

$ws.Range("B10").Value = "End User"
$ws.Range("C10").Value = "View the free videos"
$ws.Range("D10").Value = "I can select the courses"

$ws.Range("D11").Select()
